# Consolidate the various gasoline ("بنزین") product-name variants that used
# to be tracked as separate product lines into a single unified label.
#
# Previously the sheet distinguished:
#   - بنزین معمولی            (Regular gasoline)
#   - بنزین ایزومریت          (Isomerate gasoline)
#   - بنزین یورو 5            (Euro-5 gasoline)
#   - بنزین معمولی و یورو 5   (Regular & Euro-5 gasoline)
#
# These are now reported under a single product name: بنزین (Gasoline).
# This touches every yearly block on the sheet (production, sales volume,
# sales amount, cost) since each block repeats the same product rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$unified = "بنزین"

$cells = @(
    "B15", "B16", "B24", "B28",
    "B42", "B43", "B51", "B55",
    "B69", "B70", "B78", "B82",
    "B96", "B101", "B105",
    "B118", "B122", "B126",
    "B140", "B144", "B148"
)

foreach ($addr in $cells) {
    $ws.Range($addr).Value = $unified
}

# Reflect the user's on-screen state at save time: right-to-left sheet
# layout, scrolled down near the bottom block, with C18 selected.
try { $excel.ActiveWindow.DisplayRightToLeft = $true } catch {}
try { $excel.ActiveWindow.ScrollRow = 115 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
$ws.Range("C18").Select()
